$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.388.17'
$ws.Range("E2").Value = '  -0.04%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.947.42'

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.47'
$ws.Range("E5").Value = '  -0.63%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.612'
$ws.Range("E6").Value = '  -2.23%  '

# Row 7
$ws.Range("E7").Value = '  -0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.77'
$ws.Range("E8").Value = '  -3.25%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.362'
$ws.Range("E9").Value = '  -4.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0846'
$ws.Range("E10").Value = '  +3.33%  '

# Row 11
$ws.Range("E11").Value = '  +0.13%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.232.98'
$ws.Range("E12").Value = '  -1.59%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.813'
$ws.Range("E13").Value = '  -5.66%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.12'
$ws.Range("E14").Value = '  -11.53%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.49'
$ws.Range("E15").Value = '  -3.60%  '

# Row 16
$ws.Range("E16").Value = '  -5.24%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.949.82'
$ws.Range("E17").Value = '  -1.70%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.339.51'
$ws.Range("E18").Value = '  +0.05%  '

# Row 19
$ws.Range("E19").Value = '  +1.64%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.46'
$ws.Range("E20").Value = '  -1.56%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '229.07'
$ws.Range("E21").Value = '  -2.23%  '

# Row 22
$ws.Range("E22").Value = '  -5.53%  '

# Row 23
$ws.Range("E23").Value = '  +0.05%  '

# Row 24
$ws.Range("E24").Value = '  -7.62%  '

# Row 25
$ws.Range("E25").Value = '  -1.00%  '

# Row 26
$ws.Range("E26").Value = '  -9.34%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.33'
$ws.Range("E27").Value = '  -0.18%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.136'
$ws.Range("E28").Value = '  +7.57%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.24'
$ws.Range("E29").Value = '  -2.91%  '

# Row 30
$ws.Range("E30").Value = '  -1.65%  '

# Row 31
$ws.Range("E31").Value = '  -4.90%  '

# Row 32
$ws.Range("E32").Value = '  -5.85%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0636'
$ws.Range("E33").Value = '  +0.58%  '

# Row 34
$ws.Range("E34").Value = '  -3.69%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.12'
$ws.Range("E35").Value = '  -2.49%  '

# Row 36
$ws.Range("E36").Value = '  +0.11%  '

# Row 37
$ws.Range("E37").Value = '  +1.38%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.14'
$ws.Range("E38").Value = '  -5.77%  '

# Row 39
$ws.Range("E39").Value = '  -1.96%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0971'
$ws.Range("E40").Value = '  +0.93%  '

# Row 41
$ws.Range("E41").Value = '  -0.37%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.16'
$ws.Range("E42").Value = '  -6.76%  '

# Row 43
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0210'
$ws.Range("E43").Value = '  -1.37%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.352.20'
$ws.Range("E44").Value = '  -1.71%  '

# Row 45
$ws.Range("E45").Value = '  -3.97%  '

# Row 46
$ws.Range("E46").Value = '  -6.55%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.06'
$ws.Range("E47").Value = '  -5.90%  '

# Row 48
$ws.Range("E48").Value = '  -6.26%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.83'
$ws.Range("E49").Value = '  -0.41%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.26'
$ws.Range("E50").Value = '  -1.71%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.124.76'
$ws.Range("E51").Value = '  -1.82%  '
